$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Re-arrange the existing (unchanged) text blocks in rows 2-4 so
#    that the picture-filename ends up in column A and the
#    EN/ES/FR text blocks shift one column to the right (B/C/D),
#    leaving room for a new German column (E) and new numeric
#    columns (F/G). Using Cut() (not re-typing) keeps the original
#    formatting / rich-text runs intact.
# ------------------------------------------------------------------
foreach ($r in 2, 3, 4) {
    $ws.Range("F$r").Cut($ws.Range("G$r")) | Out-Null
    $ws.Range("E$r").Cut($ws.Range("F$r")) | Out-Null
    $ws.Range("A$r").Cut($ws.Range("J$r")) | Out-Null
    $ws.Range("B$r").Cut($ws.Range("A$r")) | Out-Null
    $ws.Range("J$r").Cut($ws.Range("B$r")) | Out-Null
    $ws.Range("J$r").Clear() | Out-Null
}

# ------------------------------------------------------------------
# 2) New German text (column E, rows 2-4)
# ------------------------------------------------------------------
$deInstr = @'
Anweisungen:
In dieser Aufgabe sehen Sie Bilder von linken oder rechten Händen, deren Handflächen nach oben oder unten zeigen. Die Bilder werden in verschiedenen Winkeln gedreht.
Ihre Aufgabe ist es zu bestimmen, ob das Bild einer linken oder rechten Hand entspricht.
Ihr Ziel ist es, sowohl SCHNELL als auch GENAU zu antworten.
Jedes Bild wird angezeigt, bis Sie geantwortet haben. Das nächste Bild erscheint automatisch.
'@

$deInstr2 = @'
Bitte verwenden Sie nur Ihre Zeigefinger, um zu antworten, und benutzen Sie die Tasten „S“ und „L“ auf Ihrer Tastatur:
Linke Hand = S | Rechte Hand = L
Sie müssen Ihre Hände während der gesamten Aufgabe auf der Tastatur halten.
Halten Sie Ihre Hände in derselben Position und so ruhig wie möglich.
'@

$deFeedback = @'
Nach jedem Bild erhalten Sie ein kurzes Feedback zu Ihrer Antwort:
Wenn Sie korrekt antworten, wird das entsprechende Feld grün.
Wenn Sie falsch antworten, wird das entsprechende Feld rot.
Denken Sie daran, dass Ihr Ziel darin besteht, so genau und schnell wie möglich zu antworten.
'@

$ws.Range("E2").Value = $deInstr
$ws.Range("E3").Value = $deInstr2
$ws.Range("E4").Value = $deFeedback

# give the new German cells the same "wrap text" style (s=1) used by
# the other text cells in these rows
$ws.Range("C2").Copy() | Out-Null
$ws.Range("E2:E4").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 3) Header row (row 1) - rename / relocate the field-name labels
# ------------------------------------------------------------------
$ws.Range("A1").Value = "inst_pics"
$ws.Range("B1").Value = "inst_msg_EN"
$ws.Range("C1").Value = "inst_msg_ES"
$ws.Range("D1").Value = "inst_msg_FR"
$ws.Range("E1").Value = "inst_msg_DE"
$ws.Range("F1").Value = "image_w"
$ws.Range("G1").Value = "image_h"

# ------------------------------------------------------------------
# 4) New header style: white font on a blue fill, applied to A1:G1.
#    Build it on a single cell first (so only one new font/fill/xf
#    combination is generated) then propagate with PasteSpecial.
# ------------------------------------------------------------------
$ws.Range("A1").Interior.Color = 15773696
$ws.Range("A1").Font.ThemeColor = 2

$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:G1").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 5) Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.666666666666664
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 12

# ------------------------------------------------------------------
# 6) Selection left where the author left it after the edit
# ------------------------------------------------------------------
$ws.Range("B8").Select() | Out-Null

Write-Host "done"
